# Fixed update to excel issue
# 1) Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Weekly Quantity
$ws2 = $wb.Worksheets.Item(2)   # Monthly Trend

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet after the last existing sheet
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

$data = @(
  @(45319.99999999999, 10, -36.57337856754167, 61.89533724481247),
  @(45417.99999999999, 25, -24.31819983503937, 76.89963537395286),
  @(45431.99999999999, 27, -22.150811929235, 74.48910601392286),
  @(45438.99999999999, 28, -21.83904254600911, 78.86987471978131),
  @(45473.99999999999, 33, -14.47525481199769, 87.46319161979015),
  @(45487.99999999999, 36, -17.52974015839381, 84.63008622609317),
  @(45515.99999999999, 40, -10.80574798815816, 90.42204398533332),
  @(45522.99999999999, 41, -8.55949040583101, 91.01206901511189),
  @(45543.99999999999, 44, -3.649238841750808, 93.96269200848521),
  @(45557.99999999999, 46, 0.7121688923496172, 96.32732039534361),
  @(45564.99999999999, 47, -7.452780310086221, 97.01467498627824),
  @(45571.99999999999, 48, -5.444813691147834, 99.71720718097346),
  @(45578.99999999999, 49, -1.456717977426887, 99.90779573435015),
  @(45585.99999999999, 50, -1.309246375322246, 98.64976014940767),
  @(45592.99999999999, 51, 3.012880936515485, 98.97050202541543),
  @(45599.99999999999, 52, -1.343182382533026, 102.7863173985139),
  @(45627.99999999999, 56, 0.9869600889338391, 104.8996676131508),
  @(45634.99999999999, 57, 6.264065056321156, 105.414701293428),
  @(45641.99999999999, 58, 10.03308714678861, 108.3052301790697),
  @(45648.99999999999, 59, 16.74044643475432, 109.9293996633346),
  @(45655.99999999999, 61, 12.59493249614236, 111.3474902201848),
  @(45662.99999999999, 62, 12.1188079653413, 106.3316416391734),
  @(45669.99999999999, 63, 14.56482611764048, 111.4806567265445),
  @(45676.99999999999, 64, 11.89277256224038, 117.8029551173542),
  @(45683.99999999999, 65, 15.8946196444488, 116.9849637296806),
  @(45690.99999999999, 66, 17.18413951842223, 115.3083683365439),
  @(45697.99999999999, 67, 17.00083403749505, 116.0284002533477)
)
$r = 2
foreach ($row in $data) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# 3) Copy header styling (bold, centered, bordered) from an existing header cell
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# 4) Copy the date-formatted style from column A of sheet 1 onto the ds column
$ws1.Range("A2").Copy()
$ws3.Range("A2:A28").PasteSpecial(-4122)

$excel.CutCopyMode = $false
